$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("O1").Value = "C15"
$ws.Range("O2").Select()
